$d = $word.ActiveDocument

# Title / heading (appears twice: Heading1 title and the bold run near the end)
$d.Content.Find.Execute("Play Mayan Magic Wildfire for Free: Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Mayan Magic Wildfire Free - Exciting Features & Stunning Graphics", 2)

# "What we like" bullets
$d.Content.Find.Execute("Game offers a range of wild features.", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple wild features enhance gameplay", 2)

$d.Content.Find.Execute("Attractive design, with well-executed graphics.", $true, $false, $false, $false, $false, $true, 1, $false, "Colorful graphics and immersive jungle scenery", 2)

$d.Content.Find.Execute("Mayan Wild Feature increases win multiplier by 1x.", $true, $false, $false, $false, $false, $true, 1, $false, "Maya Wilds feature with respins and sticky wilds", 2)

$d.Content.Find.Execute("Mystery Wild and Mystery Sync Features can be triggered at random.", $true, $false, $false, $false, $false, $true, 1, $false, "Mystery Wild and Mystery Sync features for additional excitement", 2)

# "What we don't like" bullets
$d.Content.Find.Execute("Mayan Magic Wildfire is a medium variance game.", $true, $false, $false, $false, $false, $true, 1, $false, "Limited number of central reels for wild triggers", 2)

$d.Content.Find.Execute("Mayan Magic Wildfire does not offer any free spins bonus.", $true, $false, $false, $false, $false, $true, 1, $false, "Medium variance may not appeal to players seeking high-risk, high-reward gameplay", 2)

# Italic meta-description run
$d.Content.Find.Execute("Read our review of Mayan Magic Wildfire, a jungle-themed online slot. Play Mayan Magic Wildfire for free and discover its gameplay features and graphic design.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the wild features and beautiful scenery in Mayan Magic Wildfire. Play now for free!", 2)
